$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Status: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (every cell that previously showed "Ready for handoff" gets the new text)
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime updates ---
$ws2.Range("K2").Value = "2016-08-15 20:46:27"
$ws3.Range("K2").Value = "2016-08-15 20:46:35"

# --- Error Detail cleared now that handback is in sync ---
$ws2.Range("P2").Value = ""
$ws3.Range("P2").Value = ""

# --- Column width updates (widen Status / zh-cn / de-de columns, narrow Error Detail) ---
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668   # Overview!E (zh-cn)
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668   # Overview!F (de-de)

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668   # zh-cn!C (Status)
$ws2.Columns.Item(16).ColumnWidth = 12.833333333333334  # zh-cn!P (Error Detail)

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668   # de-de!C (Status)
$ws3.Columns.Item(16).ColumnWidth = 12.833333333333334  # de-de!P (Error Detail)
